$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'61.570.85"
$ws.Cells.Item(2, 5).Value = '  -3.92%  '

# Row 3
$ws.Cells.Item(3, 4).Value = "'2.971.03"
$ws.Cells.Item(3, 5).Value = '  -5.12%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'540.15"
$ws.Cells.Item(5, 5).Value = '  -5.22%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'150.77"
$ws.Cells.Item(6, 5).Value = '  -6.68%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.08%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.569"
$ws.Cells.Item(8, 5).Value = '  -0.40%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'2.978.61"
$ws.Cells.Item(9, 5).Value = '  -5.26%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -2.50%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'6.12"
$ws.Cells.Item(11, 5).Value = '  -6.64%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -3.89%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'3.492.38"
$ws.Cells.Item(13, 5).Value = '  -5.05%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -2.29%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'61.628.51"
$ws.Cells.Item(15, 5).Value = '  -3.86%  '

# Row 16
$ws.Cells.Item(16, 4).Value = "'23.66"
$ws.Cells.Item(16, 5).Value = '  -5.13%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'2.971.92"
$ws.Cells.Item(17, 5).Value = '  -5.44%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'0.0000146"
$ws.Cells.Item(18, 5).Value = '  -4.17%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.68%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'12.00"
$ws.Cells.Item(20, 5).Value = '  -3.58%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'381.20"
$ws.Cells.Item(21, 5).Value = '  -4.25%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'6.66"
$ws.Cells.Item(22, 5).Value = '  -5.68%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.08%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -3.22%  '

# Row 25
$ws.Cells.Item(25, 4).Value = "'65.52"
$ws.Cells.Item(25, 5).Value = '  -3.55%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'0.470"
$ws.Cells.Item(26, 5).Value = '  -2.31%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'3.097.00"
$ws.Cells.Item(27, 5).Value = '  -5.31%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -1.87%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.48%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'0.0₃0941"
$ws.Cells.Item(30, 5).Value = '  -5.36%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'8.23"
$ws.Cells.Item(31, 5).Value = '  -5.56%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.03%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(33, 4).Value = "'1.72"
$ws.Cells.Item(33, 5).Value = '  -4.31%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(34, 4).Value = "'20.46"
$ws.Cells.Item(34, 5).Value = '  -2.72%  '

# Row 35
$ws.Cells.Item(35, 4).Value = "'160.94"
$ws.Cells.Item(35, 5).Value = '  +0.76%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'4.65"
$ws.Cells.Item(36, 5).Value = '  -2.69%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'5.92"
$ws.Cells.Item(37, 5).Value = '  -4.89%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'1.07"
$ws.Cells.Item(38, 5).Value = '  -2.24%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'1.26"
$ws.Cells.Item(39, 5).Value = '  -4.93%  '

# Row 40
$ws.Cells.Item(40, 4).Value = "'1.55"
$ws.Cells.Item(40, 5).Value = '  -6.06%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'3.91"
$ws.Cells.Item(41, 5).Value = '  -2.95%  '

# Row 42
$ws.Cells.Item(42, 4).Value = "'37.46"
$ws.Cells.Item(42, 5).Value = '  -1.96%  '

# Row 43
$ws.Cells.Item(43, 4).Value = "'2.411.39"
$ws.Cells.Item(43, 5).Value = '  -8.71%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'22.17"
$ws.Cells.Item(44, 5).Value = '  -5.47%  '

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.668"
$ws.Cells.Item(45, 5).Value = '  -2.82%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'0.0591"
$ws.Cells.Item(46, 5).Value = '  -2.91%  '

# Row 47
$ws.Cells.Item(47, 4).Value = "'5.14"
$ws.Cells.Item(47, 5).Value = '  -4.32%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.997"
$ws.Cells.Item(48, 5).Value = '  +0.07%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -2.62%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.0952"
$ws.Cells.Item(50, 5).Value = '  -1.78%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Bittensor'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(51, 4).Value = "'267.52"
$ws.Cells.Item(51, 5).Value = '  -5.77%  '
